$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grades")

$ws.Range("L11").Value = 13
$ws.Range("Q8").Value = 60

[void]$ws.Range("L11").Select()
